$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 344, shifting existing rows 344..428 down to 345..429
$ws.Rows("344:344").Insert()

# Populate the newly inserted row 344 with the new data record
$ws.Range("A344").Value = 5
$ws.Range("B344").Value = "Macroferia Regional de Talca"
$ws.Range("C344").Value = "Maule"
$ws.Range("D344").Value = 44543
$ws.Range("E344").Value = 7
$ws.Range("F344").Value = 100112004
$ws.Range("G344").Value = "Cebolla"
$ws.Range("H344").Value = "Sin especificar"
$ws.Range("I344").Value = "1a nueva(o)"
$ws.Range("J344").Value = 60000
$ws.Range("K344").Value = 1000
$ws.Range("L344").Value = 1000
$ws.Range("M344").Value = 1000
$ws.Range("N344").Value = "$/paquete 10 unidades (volumen en unidades)"
$ws.Range("O344").Value = "Región de O'Higgins"
$ws.Range("P344").Value = 100
$ws.Range("Q344").Value = 10
$ws.Range("R344").Value = "Hortaliza"
